# Natmi following Dr Hou advice
# Updates the LR-pairs table for Sema5a-Met: adds an "M2" sending/target cluster
# alongside the existing ECs/FAPs/sCs clusters, and refreshes all computed
# NATMI edge-weight statistics accordingly (4 x 4 = 16 cluster-pair rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, sending cluster (col A), target cluster (col D),
# followed by the 16 numeric values for columns E..T.
$data = @(
    @(2, "ECs", "ECs", 3,1,1.143562333333333,3.430687,0.02429958330097579,0.0242995833009758,3,1,3.167735333333333,9.503206,0.1182666224938439,0.1182666224938439,3.622502809169111,32.602525282522,0.002873829645014217,0.002873829645014218),
    @(3, "ECs", "FAPs", 3,1,1.143562333333333,3.430687,0.02429958330097579,0.0242995833009758,3,1,0.9421210000000001,2.826363,0.03517385669126484,0.03517385669126484,1.077374089042334,9.696366801381002,0.0008547100606859746,0.0008547100606859749),
    @(4, "ECs", "M2", 3,1,1.143562333333333,3.430687,0.02429958330097579,0.0242995833009758,3,1,1.442127,4.326381,0.05384145818700961,0.0538414581870096,1.649162117083,14.842459053747,0.001308324998261245,0.001308324998261245),
    @(5, "ECs", "sCs", 3,1,1.143562333333333,3.430687,0.02429958330097579,0.0242995833009758,3,1,21.232711,63.69813300000001,0.7927180626278817,0.7927180626278817,24.28092853415233,218.528356807371,0.01926271859701435,0.01926271859701436),
    @(6, "FAPs", "ECs", 3,1,40.163957,120.491871,0.8534448804146018,0.853444880414602,3,1,3.167735333333333,9.503206,0.1182666224938439,0.1182666224938439,127.2287857153807,1145.059071438426,0.1009340434912975,0.1009340434912975),
    @(7, "FAPs", "FAPs", 3,1,40.163957,120.491871,0.8534448804146018,0.853444880414602,3,1,0.9421210000000001,2.826363,0.03517385669126484,0.03517385669126484,37.83930733279701,340.553765995173,0.03001894791759686,0.03001894791759687),
    @(8, "FAPs", "M2", 3,1,40.163957,120.491871,0.8534448804146018,0.853444880414602,3,1,1.442127,4.326381,0.05384145818700961,0.0538414581870096,57.921526816539,521.293741348851,0.0459507168437602,0.0459507168437602),
    @(9, "FAPs", "sCs", 3,1,40.163957,120.491871,0.8534448804146018,0.853444880414602,3,1,21.232711,63.69813300000001,0.7927180626278817,0.7927180626278817,852.7896915974271,7675.107224376844,0.6765411721619473,0.6765411721619474),
    @(10, "M2", "ECs", 1,0.3333333333333333,0.002171666666666667,0.006515,0.00004614579680567107,0.00004614579680567108,3,1,3.167735333333333,9.503206,0.1182666224938439,0.1182666224938439,0.006879265232222223,0.06191338709,0.000005457507530493928,0.000005457507530493929),
    @(11, "M2", "FAPs", 1,0.3333333333333333,0.002171666666666667,0.006515,0.00004614579680567107,0.00004614579680567108,3,1,0.9421210000000001,2.826363,0.03517385669126484,0.03517385669126484,0.002045972771666667,0.018413754945,0.000001623125643746901,0.000001623125643746901),
    @(12, "M2", "M2", 1,0.3333333333333333,0.002171666666666667,0.006515,0.00004614579680567107,0.00004614579680567108,3,1,1.442127,4.326381,0.05384145818700961,0.0538414581870096,0.003131819135,0.028186372215,0.000002484556989218781,0.000002484556989218781),
    @(13, "M2", "sCs", 1,0.3333333333333333,0.002171666666666667,0.006515,0.00004614579680567107,0.00004614579680567108,3,1,21.232711,63.69813300000001,0.7927180626278817,0.7927180626278817,0.04611037072166668,0.414993336495,0.00003658060664221147,0.00003658060664221147),
    @(14, "sCs", "ECs", 3,1,5.751294333333334,17.253883,0.1222093904876166,0.1222093904876167,3,1,3.167735333333333,9.503206,0.1182666224938439,0.1182666224938439,18.21857827209978,163.967204448898,0.01445329185000171,0.01445329185000172),
    @(15, "sCs", "FAPs", 3,1,5.751294333333334,17.253883,0.1222093904876166,0.1222093904876167,3,1,0.9421210000000001,2.826363,0.03517385669126484,0.03517385669126484,5.418415168614334,48.76573651752901,0.004298575587338252,0.004298575587338253),
    @(16, "sCs", "M2", 3,1,5.751294333333334,17.253883,0.1222093904876166,0.1222093904876167,3,1,1.442127,4.326381,0.05384145818700961,0.0538414581870096,8.294096843047,74.646871587423,0.006579931787998942,0.006579931787998942),
    @(17, "sCs", "sCs", 3,1,5.751294333333334,17.253883,0.1222093904876166,0.1222093904876167,3,1,21.232711,63.69813300000001,0.7927180626278817,0.7927180626278817,122.1155704556044,1099.040134100439,0.09687759126227774,0.09687759126227775)

)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2  = $row[1]   # A: Sending cluster
    $ws.Cells.Item($r, 2).Value2  = "Sema5a"  # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value2  = "Met"     # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value2  = $row[2]   # D: Target cluster

    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value2 = $row[3 + $i]
    }
}
